# Updates cryptos list prices / volume(1h) columns, matching the
# "Updated cryptos list ... with GitHub Actions" commit.
# All cells in this sheet hold text (prices/percentages are formatted
# strings, not numbers), so a plain ".Value = ..." assignment would let
# Excel auto-coerce numeric-looking text ("70.923.74", "0.999", ...)
# into a real number and drop the original text formatting. Forcing the
# cell to Text format before the write keeps it a string, and resetting
# the style back to "Normal" afterwards avoids leaving a stray
# NumberFormat/quote-prefix style behind so cell styling matches the
# original workbook exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "70.923.74"
Set-TextCell 2 5 "  -3.00%  "
Set-TextCell 3 4 "3.856.24"
Set-TextCell 3 5 "  -3.42%  "
Set-TextCell 4 5 "  -0.10%  "
Set-TextCell 5 4 "598.71"
Set-TextCell 5 5 "  +0.46%  "
Set-TextCell 6 4 "169.15"
Set-TextCell 6 5 "  +3.37%  "
Set-TextCell 7 5 "  -2.51%  "
Set-TextCell 8 5 "  +0.17%  "
Set-TextCell 9 4 "0.744"
Set-TextCell 9 5 "  -0.89%  "
Set-TextCell 10 4 "0.176"
Set-TextCell 10 5 "  +3.82%  "
Set-TextCell 11 4 "53.37"
Set-TextCell 11 5 "  -2.12%  "
Set-TextCell 12 4 "0.0000322"
Set-TextCell 12 5 "  +0.45%  "
Set-TextCell 13 4 "11.36"
Set-TextCell 13 5 "  +3.33%  "
Set-TextCell 14 4 "4.456.20"
Set-TextCell 14 5 "  -3.82%  "
Set-TextCell 15 4 "21.23"
Set-TextCell 15 5 "  +3.83%  "
Set-TextCell 16 4 "3.857.28"
Set-TextCell 16 5 "  -3.37%  "
Set-TextCell 17 4 "13.91"
Set-TextCell 17 5 "  -1.73%  "
Set-TextCell 18 4 "1.21"
Set-TextCell 18 5 "  -5.02%  "
Set-TextCell 19 5 "  -2.39%  "
Set-TextCell 20 4 "70.696.36"
Set-TextCell 20 5 "  -2.95%  "
Set-TextCell 21 4 "439.44"
Set-TextCell 21 5 "  +0.34%  "
Set-TextCell 22 4 "4.77"
Set-TextCell 22 5 "  +0.62%  "
Set-TextCell 23 4 "94.66"
Set-TextCell 23 5 "  -1.79%  "
Set-TextCell 24 4 "3.29"
Set-TextCell 24 5 "  -4.76%  "
Set-TextCell 25 4 "13.85"
Set-TextCell 25 5 "  -3.56%  "
Set-TextCell 26 4 "11.58"
Set-TextCell 26 5 "  +2.11%  "
Set-TextCell 27 5 "  -8.09%  "
Set-TextCell 28 5 "  +0.15%  "
Set-TextCell 29 4 "10.51"
Set-TextCell 29 5 "  +1.21%  "
Set-TextCell 30 4 "8.55"
Set-TextCell 30 5 "  +9.55%  "
Set-TextCell 31 4 "35.01"
Set-TextCell 31 5 "  -3.68%  "
Set-TextCell 32 4 "13.53"
Set-TextCell 32 5 "  -2.10%  "
Set-TextCell 33 4 "48.39"
Set-TextCell 33 5 "  +0.26%  "
Set-TextCell 34 5 "  -4.25%  "
Set-TextCell 35 2 "PEPE"
Set-TextCell 35 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell 35 4 "0.0₃0986"
Set-TextCell 35 5 "  +9.08%  "
Set-TextCell 36 2 "OKB"
Set-TextCell 36 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 36 4 "68.81"
Set-TextCell 36 5 "  -3.00%  "
Set-TextCell 37 4 "639.24"
Set-TextCell 37 5 "  -5.13%  "
Set-TextCell 38 4 "0.434"
Set-TextCell 38 5 "  -0.73%  "
Set-TextCell 39 5 "  +0.44%  "
Set-TextCell 40 5 "  +0.07%  "
Set-TextCell 41 4 "0.999"
Set-TextCell 41 5 "  -0.33%  "
Set-TextCell 42 4 "3.24"
Set-TextCell 42 5 "  -4.35%  "
Set-TextCell 43 4 "2.88"
Set-TextCell 43 5 "  +9.36%  "
Set-TextCell 44 4 "3.16"
Set-TextCell 44 5 "  +18.29%  "
Set-TextCell 45 4 "0.0470"
Set-TextCell 45 5 "  -4.15%  "
Set-TextCell 46 4 "10.06"
Set-TextCell 46 5 "  -5.63%  "
Set-TextCell 47 4 "0.144"
Set-TextCell 47 5 "  -3.88%  "
Set-TextCell 48 4 "2.87"
Set-TextCell 48 5 "  -14.28%  "
Set-TextCell 49 4 "2.909.22"
Set-TextCell 49 5 "  +0.33%  "
Set-TextCell 50 4 "3.31"
Set-TextCell 50 5 "  -3.11%  "
Set-TextCell 51 4 "0.000277"
Set-TextCell 51 5 "  +2.57%  "
